$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
        if ($sh.Type -eq 6) {
            $found = Find-ShapeByName $sh.GroupItems $name
            if ($found) { return $found }
        }
    }
    return $null
}

# The copyright notice lives in the "TextBox 11" shape inside "Group 12"
# on slide 1. Update the copyright year from 2023 to 2024, keeping the
# rest of the line (and its run formatting) intact.
$tb = Find-ShapeByName $s.Shapes "TextBox 11"
$tr = $tb.TextFrame.TextRange
$fullText = $tr.Text

$crIdx = $fullText.IndexOf([char]13)
if ($crIdx -ge 0) {
    $line1Len = $crIdx
} else {
    $line1Len = $fullText.Length
}

$line1 = $tr.Characters(1, $line1Len)
$line1.Text = $line1.Text -replace "2011-2023", "2011-2024"
